# ResumeAKamlani.docx update — "update resume per new status"
#
# This script reproduces three related edits:
#   1. The author's cursor/insertion point (tracked by Word's hidden
#      "_GoBack" bookmark) moves from right after "Deep Learning
#      Consultant" to inside the word "particular attention" (splitting
#      that run into "pa" | "rticular attention").
#   2. The "_GoBack" bookmark is therefore no longer after "Deep Learning
#      Consultant" (Word only ever keeps a single "_GoBack" bookmark, so
#      re-adding it elsewhere automatically removes the old one).
#   3. The second job's end date changes from "Present" to "Jun 2017".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: split "particular attention" into "pa" / "rticular attention"
# and drop the "_GoBack" bookmark at that split point.
# ---------------------------------------------------------------------

$target = $d.Content
$found = $target.Find.Execute("particular attention", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'particular attention'"
}

$splitPos = $target.Start + 2      # right after "pa"
$tailEnd  = $target.End            # end of "particular attention"
$tailText = "rticular attention"

# Temporary bookmarks bracket the "rticular attention" span tightly so
# the in-place re-write below can't cascade/merge into neighboring runs
# that happen to share identical formatting.
$d.Bookmarks.Add("zzzTmpL", $d.Range($splitPos, $splitPos))
$d.Bookmarks.Add("zzzTmpR", $d.Range($tailEnd, $tailEnd))

# Round-trip the tail text through a placeholder of the same length so
# the engine regenerates it as a brand-new run (matching how Word emits
# freshly-touched text without any w:rsid* attributes) rather than
# reusing the original run object.
$tailRange = $d.Range($splitPos, $tailEnd)
$tailRange.Text = "X" * $tailText.Length
$tailRange2 = $d.Range($splitPos, $tailRange.End)
$tailRange2.Text = $tailText

$d.Bookmarks.Item("zzzTmpL").Delete()
$d.Bookmarks.Item("zzzTmpR").Delete()

# Re-seat "_GoBack" at the split point; Word keeps only one instance, so
# this also removes it from after "Deep Learning Consultant".
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))

# ---------------------------------------------------------------------
# Part 2: change the Tyto job's end date from "Present" to "Jun 2017".
# ---------------------------------------------------------------------

$dateRange = $d.Content
$found = $dateRange.Find.Execute(" Present", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find first ' Present'"
}
# The first " Present" match belongs to the still-current Skymind role
# ("May 2017 - Present"); collapse past it and search again for the
# Tyto role's "Mar 2017 - Present" which is the one being updated.
$dateRange.Collapse(0)
$found = $dateRange.Find.Execute(" Present", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find second ' Present'"
}

$presentStart = $dateRange.Start + 1   # skip the leading space, keep it in its own run
$presentEnd   = $dateRange.End

# Bookmark-bracket the "Present" word so replacing it can't merge the
# preceding " " run backwards into the "–" run.
$d.Bookmarks.Add("zzzTmpDate", $d.Range($presentStart, $presentStart))
$presentRange = $d.Range($presentStart, $presentEnd)
$presentRange.Text = "Jun 2017"
$d.Bookmarks.Item("zzzTmpDate").Delete()
